$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.000.93"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.118.33"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.47%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "347.60"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.49%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5200"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.93%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4458"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +1.13%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "54.18"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +3.80%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.09363"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +0.90%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.179"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.04%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "25.29"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.35%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "8.482"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +4.29%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.118.69"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.882"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +2.25%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "102.47"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +3.25%  "

$ws.Range("E17").Value = "  -0.18%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.37%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "21.66"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +5.00%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.06682"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.31%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.302"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").Value = "30.037.27"
$ws.Range("E23").Value = "  -0.12%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "12.78"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +1.57%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.330"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "2.381.34"
$ws.Range("E26").Value = "  +1.48%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "22.13"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.79%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.549"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -0.06%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "162.84"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +0.34%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "134.18"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +0.74%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.157"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -0.78%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.786"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +9.50%  "

$ws.Range("E33").Value = "  -0.27%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.251"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.970"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.542"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +5.27%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "10.78"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +6.87%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02611"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +2.21%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.06908"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.91%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "12.73"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.02%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.7045"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.2249"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.334"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +2.64%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.6851"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +3.35%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "14.57"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +2.52%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.357"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +3.68%  "

$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("E48").Value = "  +1.73%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.642"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.50%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.254"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +7.78%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.226"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.50%  "
